$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look like numbers stay as plain text,
# matching the source inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.029.00'
$ws.Range('E2').Value = '  -2.29%  '
$ws.Range('D3').Value = '1.830.37'
$ws.Range('E3').Value = '  -1.21%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D5').Value = '326.83'
$ws.Range('E5').Value = '  -2.62%  '
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '0.4619'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('D8').Value = '0.3868'
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').Value = '0.07866'
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').Value = '0.9585'
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('D11').Value = '21.89'
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').Value = '1.965.81'
$ws.Range('E12').Value = '  +5.39%  '
$ws.Range('D13').Value = '5.657'
$ws.Range('E13').Value = '  -3.13%  '
$ws.Range('D14').Value = '6.887'
$ws.Range('E14').Value = '  -1.87%  '
$ws.Range('D15').Value = '0.06763'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = '86.75'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '0.000009935'
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('D19').Value = '16.65'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').Value = '28.057.71'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').Value = '5.306'
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('D23').Value = '11.01'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('D24').Value = '2.099'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '2.095.22'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').Value = '153.63'
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').Value = '19.19'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').Value = '5.732'
$ws.Range('E28').Value = '  -8.45%  '
$ws.Range('D29').Value = '1.971'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('D30').Value = '117.25'
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('D31').Value = '0.9366'
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('D33').Value = '5.295'
$ws.Range('E33').Value = '  -1.81%  '
$ws.Range('D34').Value = '1.315'
$ws.Range('E34').Value = '  -2.70%  '
$ws.Range('D35').Value = '3.321'
$ws.Range('E35').Value = '  -5.06%  '
$ws.Range('D36').Value = '0.05868'
$ws.Range('E36').Value = '  -4.63%  '
$ws.Range('D37').Value = '0.02144'
$ws.Range('E37').Value = '  -2.62%  '
$ws.Range('D38').Value = '1.143'
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = '7.733'
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D41').Value = '9.891'
$ws.Range('E41').Value = '  -2.27%  '
$ws.Range('D42').Value = '0.1760'
$ws.Range('E42').Value = '  -1.66%  '
$ws.Range('D43').Value = '1.229'
$ws.Range('E43').Value = '  -1.95%  '
$ws.Range('D44').Value = '11.59'
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('D45').Value = '0.5262'
$ws.Range('E45').Value = '  -2.86%  '
$ws.Range('D46').Value = '0.07024'
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').Value = '2.150'
$ws.Range('E47').Value = '  -8.54%  '
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('D49').Value = '112.91'
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').Value = '2.322'
$ws.Range('E51').Value = '  -0.24%  '
